# Step 6 TOC entry: "Events" -> "Events" + " and Scoring" (new run),
# and relocate the singleton "_GoBack" bookmark from the end of the
# "Finish: Placement:" paragraph to the end of this (now longer) entry.

$d = $word.ActiveDocument

# 1) Find the "Events" list entry and append " and Scoring" right after it.
#    (InsertAfter on a collapsed range merges into the same run when the
#    formatting matches, so for now this lands as a single run of
#    "Events and Scoring" - that gets split back apart in step 4 below.)
$r = $d.Content
$found = $r.Find.Execute("Events", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Events' paragraph text"
}
$r.Collapse(0)
$r.InsertAfter(" and Scoring")

# Remember the character offset right after " and Scoring" - that's where
# the relocated _GoBack bookmark needs to sit.
$bmPos = $r.End

# 2) Insert a one-character placeholder right at that spot. This keeps the
#    eventual bookmark location from being the literal last character of
#    the paragraph while we add it (adding a bookmark exactly at the very
#    end of a paragraph's text mis-places it at the paragraph start), then
#    we strip the placeholder back out again afterwards.
$r.Collapse(0)
$r.InsertAfter("X")

# 3) Add the bookmark at the remembered position (now safely mid-paragraph
#    because of the trailing "X"). Bookmarks named "_GoBack" are a Word
#    singleton, so this automatically removes the old "_GoBack" bookmark
#    wherever it used to be (after "Finish: Placement:").
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 4) Remove the "X" placeholder again.
$placeholderRange = $d.Range($bmPos, $bmPos + 1)
$placeholderRange.Delete()

# 5) Re-split "Events" and " and Scoring" into two distinct runs (matching
#    the target markup) by nudging the Bold property of just the new text
#    off and back on - toggling formatting forces the engine to break it
#    into its own run instead of silently re-merging with "Events".
$r2 = $d.Content
$found2 = $r2.Find.Execute(" and Scoring", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find inserted ' and Scoring' text"
}
$r2.Font.Bold = $false
$r2.Font.Bold = $true
